# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Mapping of cell -> new value (applies identically to both sheets).
$updates = @{
    "F3"  = 3259
    "F4"  = 238
    "F5"  = 133
    "F7"  = 1707
    "F8"  = 1642
    "F10" = 377
    "F12" = 30
    "F13" = 193
    "F23" = 387
    "F24" = 232
    "F25" = 107
    "F26" = 40
    "F27" = 12
    "F29" = 337
    "F30" = 2233
    "F31" = 11
    "F33" = 474
    "F34" = 446
    "F36" = 428
    "F37" = 230
    "F38" = 351
    "F40" = 530
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
